$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("ES1:EV1").EntireColumn.Insert()
$c = $ws.Range("ES1")
$c.NumberFormat = "General"
Write-Output "done"
